$d = $word.ActiveDocument
$r = $d.Range(10, 20)
Write-Host "before:" $r.Start $r.End
$r.Collapse(1)
Write-Host "after:" $r.Start $r.End
